# Updates cryptocurrency price/volume figures in columns D (Price) and E (Volume 1h)
# on Sheet1, per the latest scrape (cryptos list refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.611.55'
$ws.Range("D3").Value = '1.597.14'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'211.28"
$ws.Range("E5").Value = '  -0.06%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("E8").Value = '  +0.25%  '
$ws.Range("D9").Value = "'0.247"
$ws.Range("E9").Value = '  +0.29%  '
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("D11").Value = "'0.0841"
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("D12").Value = '1.820.46'
$ws.Range("E12").Value = '  +0.14%  '
$ws.Range("D13").Value = '1.601.91'
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").Value = "'4.03"
$ws.Range("E14").Value = '  -0.02%  '
$ws.Range("D15").Value = "'0.523"
$ws.Range("E15").Value = '  +0.03%  '
$ws.Range("E16").Value = '  +0.27%  '
$ws.Range("D17").Value = '26.600.44'
$ws.Range("E17").Value = '  -0.15%  '
$ws.Range("E18").Value = '  +1.30%  '
$ws.Range("E19").Value = '  +0.02%  '
$ws.Range("D20").Value = "'208.17"
$ws.Range("E20").Value = '  -0.43%  '
$ws.Range("E21").Value = '  +5.10%  '
$ws.Range("E22").Value = '  +0.70%  '
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").Value = "'145.44"
$ws.Range("E25").Value = '  -0.82%  '
$ws.Range("E26").Value = '  -0.01%  '
$ws.Range("D27").Value = "'7.11"
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("E28").Value = '  -0.22%  '
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("E30").Value = '  +1.56%  '
$ws.Range("E31").Value = '  -0.10%  '
$ws.Range("E32").Value = '  +0.52%  '
$ws.Range("E33").Value = '  +1.10%  '
$ws.Range("D34").Value = '1.284.92'
$ws.Range("E34").Value = '  +0.05%  '
$ws.Range("D35").Value = "'0.617"
$ws.Range("E35").Value = '  -6.60%  '
$ws.Range("D36").Value = "'2.46"
$ws.Range("E36").Value = '  +0.84%  '
$ws.Range("E37").Value = '  +0.74%  '
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("E39").Value = '  +0.80%  '
$ws.Range("D40").Value = "'1.06"
$ws.Range("E40").Value = '  +20.94%  '
$ws.Range("D41").Value = "'5.47"
$ws.Range("E41").Value = '  +2.06%  '
$ws.Range("D42").Value = "'2.19"
$ws.Range("E42").Value = '  -0.25%  '
$ws.Range("D43").Value = "'64.32"
$ws.Range("E43").Value = '  +1.27%  '
$ws.Range("E44").Value = '  -1.06%  '
$ws.Range("D45").Value = '1.733.19'
$ws.Range("E45").Value = '  +0.18%  '
$ws.Range("D46").Value = "'90.08"
$ws.Range("E46").Value = '  +0.42%  '
$ws.Range("D47").Value = "'1.61"
$ws.Range("E47").Value = '  -1.28%  '
$ws.Range("E49").Value = '  +1.04%  '
$ws.Range("E50").Value = '  -0.14%  '
$ws.Range("D51").Value = "'7.44"
$ws.Range("E51").Value = '  -0.91%  '
